## Fruta / hortaliza, semanal
## Insert the new week's two rows (Especial / Primera) at the top of the
## "Macroferia Regional de Talca - Chirimoya" data block (row 85) and push
## everything else down by two rows. The two oldest data rows that fall off
## the bottom of the historical block simply reappear at the new end of the
## range (rows 127-128) with their original values, which is exactly what a
## native row-insert gives us for free.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 85; rows 85-126 (and their
# formatting) shift down to 87-128 automatically.
$ws.Rows("85:86").Insert()

# Shared boilerplate values that are identical for every data row in this
# subset (only D, L, M, N, O, P, S vary row to row).
$colA = 5
$colB = "Macroferia Regional de Talca"
$colC = "Maule"
$colE = 7
$colF = "Fruta"
$colG = 100107
$colH = "Otros"
$colI = 100107002
$colJ = "Chirimoya"
$colK = "Cultivar IV Región"
$colQ = "`$/bandeja 10 kilos"
$colR = "Provincia de Limarí"
$colT = 10

function Set-ChirimoyaRow($r, $fecha, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $precioKg) {
    $ws.Cells.Item($r, 1).Value = $colA
    $ws.Cells.Item($r, 2).Value = $colB
    $ws.Cells.Item($r, 3).Value = $colC
    $ws.Cells.Item($r, 4).Value = $fecha
    $ws.Cells.Item($r, 5).Value = $colE
    $ws.Cells.Item($r, 6).Value = $colF
    $ws.Cells.Item($r, 7).Value = $colG
    $ws.Cells.Item($r, 8).Value = $colH
    $ws.Cells.Item($r, 9).Value = $colI
    $ws.Cells.Item($r, 10).Value = $colJ
    $ws.Cells.Item($r, 11).Value = $colK
    $ws.Cells.Item($r, 12).Value = $calidad
    $ws.Cells.Item($r, 13).Value = $volumen
    $ws.Cells.Item($r, 14).Value = $precioMin
    $ws.Cells.Item($r, 15).Value = $precioMax
    $ws.Cells.Item($r, 16).Value = $precioProm
    $ws.Cells.Item($r, 17).Value = $colQ
    $ws.Cells.Item($r, 18).Value = $colR
    $ws.Cells.Item($r, 19).Value = $precioKg
    $ws.Cells.Item($r, 20).Value = $colT
}

# New week's data (2022-11-10 / serial 44875).
Set-ChirimoyaRow 85 44875 "Especial" 200 25000 25000 25000 2500
Set-ChirimoyaRow 86 44875 "Primera"  250 22000 22000 22000 2200

"OK"
